$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the confidential note's date from 2021-04-28 to 2021-04-29
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-29 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.4863055220826394
$ws.Range("E2").Value = -0.002320185614849146

$ws.Range("D3").Value = 0.3320276834075236
$ws.Range("E3").Value = 0.01094749079635737

$ws.Range("D4").Value = 0.09677068151956866
$ws.Range("E4").Value = 0.00877035976373719

$ws.Range("D5").Value = 0.05405291717494637
$ws.Range("E5").Value = 0.00114442664225245

$ws.Range("D6").Value = 0.03084319581532188
$ws.Range("E6").Value = 0.0126353790613718

$ws.Range("E7").Value = 0.003806839692087616
